$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = 346.5966732226262
$ws.Range("D1").Value = 1
$ws.Range("E1").Value = 0.02129094035623208
$ws.Range("F1").Value = 29.26509169136127
$ws.Range("G1").Value = 20.73490830863873
$ws.Range("H1").Value = 89.4

$ws.Range("C2").Value = 110.8307535007674
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 0.01934624761484914
$ws.Range("F2").Value = 24.22316825551013
$ws.Range("G2").Value = 19.77683174448987
$ws.Range("H2").Value = 70
$ws.Range("J2").Value = 25.3

$ws.Range("C3").Value = 2000
$ws.Range("D3").Value = 0.9
$ws.Range("E3").Value = 0.0118544253095994
$ws.Range("F3").Value = 26.64271023903579
$ws.Range("G3").Value = 33.35728976096421
$ws.Range("H3").Value = 89.4

$ws.Range("C4").Value = 2000
$ws.Range("D4").Value = 0.9
$ws.Range("E4").Value = 0.01186158014437543
$ws.Range("F4").Value = 26.74583139691585
$ws.Range("G4").Value = 33.25416860308415
$ws.Range("H4").Value = 89.4

$ws.Range("C5").Value = 2000
$ws.Range("D5").Value = 0.9
$ws.Range("E5").Value = 0.01186518096240248
$ws.Range("F5").Value = 26.799048324259
$ws.Range("G5").Value = 33.200951675741
$ws.Range("H5").Value = 89.4

$ws.Range("C6").Value = 2000
$ws.Range("D6").Value = 0.9
$ws.Range("E6").Value = 0.01186652409104374
$ws.Range("F6").Value = 26.81913973301276
$ws.Range("G6").Value = 33.18086026698724
$ws.Range("H6").Value = 89.4

$ws.Range("C7").Value = 1126.140135161729
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0.01787332804078674
$ws.Range("F7").Value = 21.83960164549232
$ws.Range("G7").Value = 24.16039835450768
$ws.Range("H7").Value = 89.4

$ws.Range("C8").Value = 377.0495732629889
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 0.02142240840060318
$ws.Range("F8").Value = 28.93271488473224
$ws.Range("G8").Value = 21.06728511526776
$ws.Range("H8").Value = 89.4

$ws.Range("C9").Value = 598.7418561754057
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 0.02203339332811525
$ws.Range("F9").Value = 26.5632097333916
$ws.Range("G9").Value = 23.4367902666084
$ws.Range("H9").Value = 89.4

$ws.Range("C10").Value = 204.3455593486477
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 0.0183952282767417
$ws.Range("F10").Value = 23.60136349206766
$ws.Range("G10").Value = 18.39863650793234
$ws.Range("H10").Value = 70
$ws.Range("J10").Value = 25.3

$ws.Range("C11").Value = 887.1812137715135
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 0.02194557291520167
$ws.Range("F11").Value = 27.08673342122646
$ws.Range("G11").Value = 22.91326657877354
$ws.Range("H11").Value = 89.4

$ws.Range("C12").Value = 325.7673349568731
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 0.02077141887055474
$ws.Range("F12").Value = 30.35239132663112
$ws.Range("G12").Value = 19.64760867336888
$ws.Range("H12").Value = 89.4

$ws.Range("C13").Value = 2000
$ws.Range("D13").Value = 0.9
$ws.Range("E13").Value = 0.01186534198719947
$ws.Range("F13").Value = 26.80145000768432
$ws.Range("G13").Value = 33.19854999231568
$ws.Range("H13").Value = 89.4

$ws.Range("C14").Value = 2000
$ws.Range("D14").Value = 0.9
$ws.Range("E14").Value = 0.01185893647689549
$ws.Range("F14").Value = 26.70731588552928
$ws.Range("G14").Value = 33.29268411447072
$ws.Range("H14").Value = 89.4

$ws.Range("C15").Value = 2000
$ws.Range("D15").Value = 0.9
$ws.Range("E15").Value = 0.01186091083911719
$ws.Range("F15").Value = 26.73604083024924
$ws.Range("G15").Value = 33.26395916975076
$ws.Range("H15").Value = 89.4

$ws.Range("C16").Value = 2000
$ws.Range("D16").Value = 0.9
$ws.Range("E16").Value = 0.01186679790907952
$ws.Range("F16").Value = 26.82325215286034
$ws.Range("G16").Value = 33.17674784713967
$ws.Range("H16").Value = 89.4
